$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of trade data to append (rows 10-12)
$data = @(
    @(9956.3700000000008, 10046.790000000001, 18.84, 19.010000000000002, $true, 0.9, 42613.765555555554, $false),
    @(10011.129999999999, 9956.3700000000008, 18.93, 18.824999999999999, $true, -0.55000000000000004, 42614.672905092593, $true),
    @(10064.19, 10011.129999999999, 18.72, 18.62, $true, -0.53, 42615.750162037039, $true)
)

$startRow = 10
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]

    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $ws.Cells.Item($row, 5).Value = $rowData[4]
    $ws.Cells.Item($row, 6).Value = $rowData[5]
    $ws.Cells.Item($row, 7).Value = $rowData[6]
    $ws.Cells.Item($row, 8).Value = $rowData[7]
}
